# Reflects the Dropbox changes as of March 25th for "Tareas Mario".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare the new date cells (rows 15-20, column B) to use the same ---
# --- date format/style as the rest of the B column (copy format only)  ---
$ws.Range("B5").Copy()
$ws.Range("B15:B20").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 15: reunion with the implementation team ---
$ws.Range("A15").Value = "Reunión con el equipo de implementación, para tratar aspectos de la arquitectura de la interfaz"
$ws.Range("B15").Value = 40991
$ws.Range("C15").Value = 60

# --- Row 16: restructuring the interface architecture ---
$ws.Range("A16").Value = "Reestructuración de la arquitectura de la interfaz"
$ws.Range("B16").Value = 40991
$ws.Range("C16").Value = 240
$ws.Range("D16").Value = "Debido a que se hizo sin atender a lo dicho por el grupo de diseño por parte de otro compañero"

# --- Row 17: finishing touches on the restructuring ---
$ws.Range("A17").Value = "Ajustar ultimos detalles de la reestructuración para poder seguir añadiendo cosas"
$ws.Range("B17").Value = 40991
$ws.Range("C17").Value = 60

# --- Row 18: adding almost everything remaining to the interface ---
$ws.Range("A18").Value = "Agregar casi todo lo restante a la interfaz (quedan los Datos Personales de Voluntarios y Beneficiarios, 1 solo panel)"
$ws.Range("B18").Value = 40991
$ws.Range("C18").Value = 240

# --- Row 19: helper controller and beneficiary panel design ---
$ws.Range("A19").Value = "Creado el controlador de ayudas y el diseño del panel de datos de beneficiarios"
$ws.Range("B19").Value = 40992
$ws.Range("C19").Value = 240

# --- Row 20: finishing up the interface and testing event handlers ---
$ws.Range("A20").Value = "Puesta a punto de la interfaz y prueba de los manejadores de eventos de la misma"
$ws.Range("B20").Value = 40993
$ws.Range("C20").Value = 210

# --- Fix the accent on the duration header (C4) ---
$ws.Range("C4").Value = "DURACIÓN (min)"

# --- Update window/selection state ---
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 2
[void]$ws.Range("C5").Select()

Write-Host "Edit applied"
